$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("G5").Value = 2.4
$ws.Range("I5").Value = 2.9
$ws.Range("K5").Value = 10
$ws.Range("R5").Value = 1.8
$ws.Range("S5").Value = 1.95
$ws.Range("X5").Value = 19
$ws.Range("Z5").Value = 10

# Row 6
$ws.Range("N6").Value = 2.03
$ws.Range("O6").Value = 1.87

# Row 7
$ws.Range("G7").Value = 3.2
$ws.Range("I7").Value = 2.35
$ws.Range("T7").Value = 8
$ws.Range("U7").Value = 15
$ws.Range("AF7").Value = 10
$ws.Range("AH7").Value = 21
$ws.Range("AI7").Value = 21

# Row 12
$ws.Range("G12").Value = 2.67
$ws.Range("H12").Value = 2.77
$ws.Range("I12").Value = 2.82
$ws.Range("L12").Value = 1.5
$ws.Range("M12").Value = 2.25
$ws.Range("N12").Value = 2.42
$ws.Range("O12").Value = 1.44
$ws.Range("P12").Value = 1.53
$ws.Range("Q12").Value = 2.18
$ws.Range("R12").Value = 2
$ws.Range("S12").Value = 1.65
$ws.Range("T12").Value = 6.5
$ws.Range("U12").Value = 12
$ws.Range("V12").Value = 10.5
$ws.Range("W12").Value = 32
$ws.Range("X12").Value = 28
$ws.Range("Y12").Value = 45
$ws.Range("Z12").Value = 6.1
$ws.Range("AA12").Value = 5.6
$ws.Range("AB12").Value = 17
$ws.Range("AC12").Value = 110
$ws.Range("AE12").Value = 6.6
$ws.Range("AF12").Value = 12.5
$ws.Range("AG12").Value = 11
$ws.Range("AH12").Value = 35
$ws.Range("AI12").Value = 30
$ws.Range("AJ12").Value = 50

# Row 19
$ws.Range("G19").Value = 2.32
$ws.Range("H19").Value = 3
$ws.Range("I19").Value = 3.1
$ws.Range("R19").Value = 1.8
$ws.Range("S19").Value = 1.8
$ws.Range("T19").Value = 6.8
$ws.Range("U19").Value = 10.5
$ws.Range("V19").Value = 9.25
$ws.Range("X19").Value = 21
$ws.Range("Y19").Value = 35
$ws.Range("AA19").Value = 5.9
$ws.Range("AF19").Value = 15.5
$ws.Range("AG19").Value = 10.75
$ws.Range("AH19").Value = 40
$ws.Range("AI19").Value = 28
$ws.Range("AJ19").Value = 37

# Row 20
$ws.Range("G20").Value = 2.8
$ws.Range("H20").Value = 3.5
$ws.Range("I20").Value = 2.25
$ws.Range("S20").Value = 2.32
$ws.Range("T20").Value = 12
$ws.Range("U20").Value = 17
$ws.Range("W20").Value = 35
$ws.Range("Y20").Value = 24
$ws.Range("AA20").Value = 7.2
$ws.Range("AB20").Value = 11.5
$ws.Range("AE20").Value = 11
$ws.Range("AF20").Value = 13.5
$ws.Range("AH20").Value = 24

# Row 21
$ws.Range("G21").Value = 1.5
$ws.Range("I21").Value = 7.5
$ws.Range("J21").Value = 1.07
$ws.Range("K21").Value = 9
$ws.Range("Z21").Value = 8.5
$ws.Range("AJ21").Value = 67

# Row 22
$ws.Range("J22").Value = 1.06
$ws.Range("K22").Value = 10

# Row 23
$ws.Range("G23").Value = 3.6
$ws.Range("H23").Value = 2.9
$ws.Range("I23").Value = 2.2
$ws.Range("J23").Value = 1.08
$ws.Range("K23").Value = 8
$ws.Range("N23").Value = 2.2
$ws.Range("O23").Value = 1.65
$ws.Range("T23").Value = 9.5
$ws.Range("U23").Value = 17
$ws.Range("X23").Value = 29
$ws.Range("AA23").Value = 5.5
$ws.Range("AF23").Value = 10
$ws.Range("AG23").Value = 9.5
$ws.Range("AH23").Value = 21
$ws.Range("AJ23").Value = 34

# Row 24
$ws.Range("G24").Value = 2.7
$ws.Range("H24").Value = 3
$ws.Range("I24").Value = 2.8
$ws.Range("P24").Value = 1.3
$ws.Range("Q24").Value = 3.4
$ws.Range("W24").Value = 26
$ws.Range("X24").Value = 19
$ws.Range("Y24").Value = 23
$ws.Range("AD24").Value = 101

# Row 25
$ws.Range("G25").Value = 4
$ws.Range("I25").Value = 2
$ws.Range("J25").Value = 1.05
$ws.Range("K25").Value = 11
$ws.Range("V25").Value = 13
$ws.Range("X25").Value = 29
$ws.Range("AE25").Value = 8

# Row 26
$ws.Range("H26").Value = 3.9
$ws.Range("I26").Value = 5.25
$ws.Range("K26").Value = 12
$ws.Range("L26").Value = 1.22
$ws.Range("M26").Value = 4
$ws.Range("AD26").Value = 201

# Row 27
$ws.Range("G27").Value = 3.5
$ws.Range("H27").Value = 3.7
$ws.Range("I27").Value = 2
$ws.Range("K27").Value = 13
$ws.Range("R27").Value = 1.7
$ws.Range("S27").Value = 2.05
$ws.Range("T27").Value = 11
$ws.Range("U27").Value = 19
$ws.Range("V27").Value = 12
$ws.Range("X27").Value = 26
$ws.Range("AC27").Value = 41
$ws.Range("AF27").Value = 10

# Row 28
$ws.Range("G28").Value = 3.2
$ws.Range("I28").Value = 2.4
$ws.Range("L28").Value = 1.25
$ws.Range("M28").Value = 3.75
$ws.Range("N28").Value = 1.85
$ws.Range("O28").Value = 1.95
$ws.Range("T28").Value = 11
$ws.Range("U28").Value = 17
$ws.Range("X28").Value = 23
$ws.Range("Z28").Value = 10
$ws.Range("AE28").Value = 9.5
$ws.Range("AH28").Value = 23
$ws.Range("AI28").Value = 19
$ws.Range("AJ28").Value = 26

# Row 29
$ws.Range("G29").Value = 1.42
$ws.Range("H29").Value = 4.75
$ws.Range("J29").Value = 1.03
$ws.Range("K29").Value = 15
$ws.Range("L29").Value = 1.18
$ws.Range("M29").Value = 4.5
$ws.Range("N29").Value = 1.62
$ws.Range("O29").Value = 2.25
$ws.Range("P29").Value = 1.29
$ws.Range("Q29").Value = 3.5
$ws.Range("T29").Value = 8
$ws.Range("W29").Value = 9.5
$ws.Range("X29").Value = 11
$ws.Range("Z29").Value = 15
$ws.Range("AA29").Value = 9
$ws.Range("AE29").Value = 21
$ws.Range("AF29").Value = 41

# Row 30
$ws.Range("G30").Value = 2.88
$ws.Range("H30").Value = 3.3
$ws.Range("I30").Value = 2.38
$ws.Range("J30").Value = 1.07
$ws.Range("K30").Value = 9
$ws.Range("T30").Value = 8
$ws.Range("U30").Value = 13
$ws.Range("V30").Value = 11
$ws.Range("W30").Value = 29
$ws.Range("X30").Value = 26
$ws.Range("Y30").Value = 34
$ws.Range("AE30").Value = 7
$ws.Range("AF30").Value = 11
$ws.Range("AG30").Value = 10
$ws.Range("AH30").Value = 23
$ws.Range("AI30").Value = 21
$ws.Range("AJ30").Value = 34

# Row 34
$ws.Range("G34").Value = 4
$ws.Range("I34").Value = 1.95
$ws.Range("K34").Value = 9
$ws.Range("T34").Value = 10
$ws.Range("U34").Value = 19
$ws.Range("V34").Value = 13
$ws.Range("AH34").Value = 17

# Row 37
$ws.Range("G37").Value = 2.38
$ws.Range("I37").Value = 2.8
$ws.Range("X37").Value = 17

# Row 39
$ws.Range("G39").Value = 1.62
$ws.Range("H39").Value = 4.25
$ws.Range("I39").Value = 4.4
$ws.Range("L39").Value = 1.16
$ws.Range("M39").Value = 4.5
$ws.Range("O39").Value = 2.37
$ws.Range("P39").Value = 1.27
$ws.Range("Q39").Value = 3.4
$ws.Range("S39").Value = 2.25
$ws.Range("U39").Value = 9.5
$ws.Range("W39").Value = 13
$ws.Range("X39").Value = 11.5
$ws.Range("AA39").Value = 8.5
$ws.Range("AB39").Value = 13.5
$ws.Range("AE39").Value = 17.5
$ws.Range("AF39").Value = 29
$ws.Range("AG39").Value = 14.5
$ws.Range("AH39").Value = 70
$ws.Range("AI39").Value = 35
$ws.Range("AJ39").Value = 35

# Row 40
$ws.Range("G40").Value = 9.25
$ws.Range("I40").Value = 1.26
$ws.Range("O40").Value = 2.67
$ws.Range("Q40").Value = 3.7
$ws.Range("R40").Value = 1.83
$ws.Range("S40").Value = 1.88
$ws.Range("T40").Value = 30
$ws.Range("U40").Value = 70
$ws.Range("V40").Value = 28
$ws.Range("X40").Value = 100
$ws.Range("Y40").Value = 75
$ws.Range("AE40").Value = 9.5
$ws.Range("AF40").Value = 7.3
$ws.Range("AJ40").Value = 24

# Row 41
$ws.Range("K41").Value = 8.75
$ws.Range("L41").Value = 1.19
$ws.Range("M41").Value = 4.15
$ws.Range("O41").Value = 2.22
$ws.Range("P41").Value = 1.31
$ws.Range("Q41").Value = 3.15
$ws.Range("T41").Value = 11.75
$ws.Range("Y41").Value = 23
$ws.Range("Z41").Value = 8.75
$ws.Range("AA41").Value = 7.1
$ws.Range("AB41").Value = 11
